$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 2.14
$ws.Range("Y2").Value = 22
$ws.Range("U3").Value = 1.48
$ws.Range("V4").Value = 3.2
